$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 80

# Column A holds a date-like text ("2025/10/08"). A direct .Value assignment of a
# date-shaped string gets auto-parsed into a serial date number (with date
# NumberFormat) by the COM layer, which does not match the source data - the
# existing rows store it as literal text. Routing it through a text formula +
# copy/paste-values round trip keeps it as plain text without Excel's
# autodetection kicking in and without minting a stray number-format style.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="2025/10/08"'
$scratch.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)
$scratch.ClearContents()

$ws.Range("B" + $newRow).Value = "水"
$ws.Range("C" + $newRow).Value = 21
$ws.Range("D" + $newRow).Value = 134
